# Daily attendance processing - normalize the "Recorded By" (column G) values.
# The recorder list for System-augmented rows should read with "System" moved
# ahead of the lowercase duplicate / the single human recorder, e.g.:
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System"     -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    } elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
